# Add 2022-Q1 data:
#  - Rename the existing "总计" sheet to "2022-Q1" and replace its contents
#    with the fund-holding breakdown for 2022-Q1.
#  - Insert a fresh "总计" sheet (after "2022-Q1") with the historical
#    quarter-over-quarter summary, including the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Create the replacement "总计" summary sheet right after it.
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: populate "2022-Q1" with the fund holding table
# ---------------------------------------------------------------------

# Header row: extend the existing bold/bordered header style (currently
# only on B1:D1) across the new E1:H1 columns, then fill in the text.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("A1:H1").Application.CutCopyMode = $false

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows (A2:H15): column A mirrors the numbered-row style already
# used by the sheet (s="2"); columns B, D, E, F, G must stay TEXT even
# though they look numeric (fund codes with leading zeros, formatted
# decimal strings) - force that via a temporary "@" number format.
$q1.Range("A2").Copy()
$q1.Range("A3:A15").PasteSpecial(-4122)
$q1.Range("A1:A15").Application.CutCopyMode = $false

$q1.Range("B2:B15").NumberFormat = "@"
$q1.Range("D2:G15").NumberFormat = "@"

$rows = @(
    @{ n="0";  code="001304"; name="建信鑫安回报灵活配置混合";           scale="2.13"; pos="66.83"; pct="5.87"; mv="0.1250"; rank=5  },
    @{ n="1";  code="002585"; name="建信兴利灵活配置混合";               scale="2.04"; pos="61.22"; pct="5.21"; mv="0.1063"; rank=5  },
    @{ n="2";  code="005005"; name="中金金泽量化精选混合A";             scale="1.69"; pos="75.38"; pct="5.92"; mv="0.1000"; rank=6  },
    @{ n="3";  code="003831"; name="建信鑫瑞回报灵活配置混合";           scale="1.80"; pos="70.33"; pct="5.44"; mv="0.0979"; rank=4  },
    @{ n="4";  code="013659"; name="中融金融鑫选3个月持有混合A";         scale="2.83"; pos="72.19"; pct="3.18"; mv="0.0900"; rank=9  },
    @{ n="5";  code="165310"; name="建信沪深300指数增强（LOF）A";       scale="2.98"; pos="92.37"; pct="2.13"; mv="0.0635"; rank=7  },
    @{ n="6";  code="000877"; name="华泰柏瑞量化优选灵活配置混合";       scale="5.61"; pos="84.68"; pct="1.12"; mv="0.0628"; rank=9  },
    @{ n="7";  code="001074"; name="华泰柏瑞量化驱动灵活配置混合A";     scale="3.76"; pos="89.85"; pct="0.98"; mv="0.0368"; rank=10 },
    @{ n="8";  code="013660"; name="中融金融鑫选3个月持有混合C";         scale="1.15"; pos="72.19"; pct="3.18"; mv="0.0366"; rank=9  },
    @{ n="9";  code="005055"; name="华泰柏瑞量化阿尔法灵活配置混合A";   scale="2.53"; pos="89.49"; pct="1.05"; mv="0.0266"; rank=7  },
    @{ n="10"; code="009208"; name="建信沪深300指数增强（LOF）C";       scale="0.30"; pos="92.37"; pct="2.13"; mv="0.0064"; rank=7  },
    @{ n="11"; code="006531"; name="华泰柏瑞量化驱动灵活配置混合C";     scale="0.08"; pos="89.85"; pct="0.98"; mv="0.0008"; rank=10 },
    @{ n="12"; code="005006"; name="中金金泽量化精选混合C";             scale="0.01"; pos="75.38"; pct="5.92"; mv="0.0006"; rank=6  },
    @{ n="13"; code="006532"; name="华泰柏瑞量化阿尔法灵活配置混合C";   scale="0.01"; pos="89.49"; pct="1.05"; mv="0.0001"; rank=7  }
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = [int]$row.n
    $q1.Cells.Item($r, 2).Value = $row.code
    $q1.Cells.Item($r, 3).Value = $row.name
    $q1.Cells.Item($r, 4).Value = $row.scale
    $q1.Cells.Item($r, 5).Value = $row.pos
    $q1.Cells.Item($r, 6).Value = $row.pct
    $q1.Cells.Item($r, 7).Value = $row.mv
    $q1.Cells.Item($r, 8).Value = [int]$row.rank
    $r = $r + 1
}

$q1.Range("B2:B15").ClearFormats()
$q1.Range("D2:G15").ClearFormats()

# ---------------------------------------------------------------------
# Step 3: populate the new "总计" sheet with the quarterly summary
# ---------------------------------------------------------------------

# Reuse the exact header/row-number styles already present on "2022-Q1"
# (s="2" in the original file) instead of rebuilding them from scratch,
# so no new style entries are introduced.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("A1:D1").Application.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$summary = @(
    @{ n=0; date="2022-Q1"; cnt=14; mv=0.75 },
    @{ n=1; date="2021-Q4"; cnt=23; mv=3.79 },
    @{ n=2; date="2021-Q3"; cnt=14; mv=1.72 },
    @{ n=3; date="2021-Q2"; cnt=1;  mv=0.08 },
    @{ n=4; date="2021-Q1"; cnt=15; mv=1.58 },
    @{ n=5; date="2020-Q4"; cnt=9;  mv=0.33 }
)

$r = 2
foreach ($row in $summary) {
    $total.Cells.Item($r, 1).Value = $row.n
    $total.Cells.Item($r, 2).Value = $row.date
    $total.Cells.Item($r, 3).Value = $row.cnt
    $total.Cells.Item($r, 4).Value = $row.mv
    $r = $r + 1
}

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$total.Range("A2:A7").Application.CutCopyMode = $false
for ($i = 0; $i -lt $summary.Count; $i++) {
    $total.Cells.Item($i + 2, 1).Value = $summary[$i].n
}

# Adding the new sheet shifts the active tab onto it; restore the
# original selection (first sheet) since the diff doesn't touch bookViews.
$wb.Worksheets.Item(1).Select()
